$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.962.52"
$ws.Range("E2").Value = "  +2.30%  "
$ws.Range("D3").Value = "1.597.20"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'210.97"
$ws.Range("E5").Value = "  +2.06%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "'0.482"
$ws.Range("E7").Value = "  +0.87%  "
$ws.Range("D8").Value = "'0.246"
$ws.Range("E8").Value = "  +0.79%  "
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("D10").Value = "'18.04"
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("D11").Value = "'0.0810"
$ws.Range("E11").Value = "  +3.51%  "
$ws.Range("D12").Value = "1.821.75"
$ws.Range("E12").Value = "  +2.29%  "
$ws.Range("D13").Value = "1.599.57"
$ws.Range("E13").Value = "  +2.37%  "
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").Value = "'0.514"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").Value = "25.978.13"
$ws.Range("E16").Value = "  +2.47%  "
$ws.Range("D17").Value = "'60.08"
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("D18").Value = "0.0₃0722"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "'200.22"
$ws.Range("E20").Value = "  +6.95%  "
$ws.Range("E21").Value = "  +2.23%  "
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("E23").Value = "  +2.19%  "
$ws.Range("E24").Value = "  +7.03%  "
$ws.Range("D25").Value = "'141.94"
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -7.69%  "
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("D29").Value = "'6.46"
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("E30").Value = "  +1.50%  "
$ws.Range("D31").Value = "'0.0476"
$ws.Range("E31").Value = "  +1.92%  "
$ws.Range("D32").Value = "'3.10"
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("D33").Value = "'2.96"
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("D34").Value = "'1.48"
$ws.Range("E34").Value = "  -0.71%  "
$ws.Range("E35").Value = "  +2.64%  "
$ws.Range("D36").Value = "1.124.78"
$ws.Range("E36").Value = "  +3.22%  "
$ws.Range("E37").Value = "  +9.26%  "
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("D40").Value = "'0.787"
$ws.Range("E40").Value = "  +1.19%  "
$ws.Range("D41").Value = "'0.490"
$ws.Range("E41").Value = "  -1.52%  "
$ws.Range("D42").Value = "'0.783"
$ws.Range("E42").Value = "  -1.79%  "
$ws.Range("D43").Value = "1.732.33"
$ws.Range("E43").Value = "  +2.06%  "
$ws.Range("D44").Value = "'5.11"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("D45").Value = "'92.77"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("E46").Value = "  +0.89%  "
$ws.Range("D47").Value = "'53.42"
$ws.Range("E47").Value = "  +1.39%  "
$ws.Range("D48").Value = "'0.0504"
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("D51").Value = "'7.18"
$ws.Range("E51").Value = "  -0.75%  "
